# Apply updated cryptocurrency price/volume data (and two row re-orderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.723.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.733.17'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.734.23'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("E10").Value = '  +2.62%  '
$ws.Range("E11").Value = '  +3.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.361.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.734.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.758.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.78%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '495.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.26%  '
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("E25").Value = '  -5.16%  '
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.881.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.668.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '433.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.82%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("B47").Value = 'Arweave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.32%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.742.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.26%  '
